$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9165092706680298
$ws.Range("B1").Value = 1.84540593624115
$ws.Range("C1").Value = 3.413103818893433
$ws.Range("D1").Value = 3.789582014083862
$ws.Range("E1").Value = 0.9568474888801575
